$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.0147
$ws.Range("G2").Value = -0.3469026548672566
$ws.Range("H2").Value = -0.3469026548672566
$ws.Range("I2").Value = -0.5530973451327433
$ws.Range("J2").Value = -0.5530973451327433
$ws.Range("K2").Value = -12.9
$ws.Range("L2").Value = -0.5707964601769911
$ws.Range("U2").Value = 29.7
$ws.Range("V2").Value = 0.1229304635761589
$ws.Range("W2").Value = -0.8543046357615894
$ws.Range("X2").Value = 0.06519206562437797
$ws.Range("Y2").Value = -0.9194967013859674
$ws.Range("Z2").Value = 0.8971814211988886
$ws.Range("AA2").Value = -0.4962286621675269
$ws.Range("AB2").Value = 0.06299505090523511
$ws.Range("AC2").Value = -0.559223713072762
$ws.Range("AD2").Value = 14.2
$ws.Range("AF2").Value = 14.2
$ws.Range("AG2").Value = -15.5
$ws.Range("AH2").Value = 0.05551211884284597
$ws.Range("AI2").Value = 0.3264367816091954
$ws.Range("AJ2").Value = -0.06855373728438745
$ws.Range("AK2").Value = -1.123188405797101
$ws.Range("AL2").Value = 1.32
$ws.Range("AM2").Value = 1.32
$ws.Range("AN2").Value = -1.339622641509434
$ws.Range("AO2").Value = -9.469696969696969
$ws.Range("AP2").Value = 1.462264150943396
$ws.Range("AQ2").Value = -9.469696969696969
$ws.Range("D3").Value = -0.0147
$ws.Range("G3").Value = -0.3469026548672566
$ws.Range("H3").Value = -0.3469026548672566
$ws.Range("I3").Value = -0.5530973451327433
$ws.Range("J3").Value = -0.5530973451327433
$ws.Range("K3").Value = -12.9
$ws.Range("L3").Value = -0.5707964601769911
$ws.Range("U3").Value = 29.7
$ws.Range("V3").Value = 0.1229304635761589
$ws.Range("W3").Value = -0.8543046357615894
$ws.Range("X3").Value = 0.06519206562437797
$ws.Range("Y3").Value = -0.9194967013859674
$ws.Range("Z3").Value = 0.8971814211988886
$ws.Range("AA3").Value = -0.4962286621675269
$ws.Range("AB3").Value = 0.06299505090523511
$ws.Range("AC3").Value = -0.559223713072762
$ws.Range("AD3").Value = 14.2
$ws.Range("AF3").Value = 14.2
$ws.Range("AG3").Value = -15.5
$ws.Range("AH3").Value = 0.05551211884284597
$ws.Range("AI3").Value = 0.3264367816091954
$ws.Range("AJ3").Value = -0.06855373728438745
$ws.Range("AK3").Value = -1.123188405797101
$ws.Range("AL3").Value = 1.32
$ws.Range("AM3").Value = 1.32
$ws.Range("AN3").Value = -1.339622641509434
$ws.Range("AO3").Value = -9.469696969696969
$ws.Range("AP3").Value = 1.462264150943396
$ws.Range("AQ3").Value = -9.469696969696969
